$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2662.25
$ws.Range("I6").Value = 177
$ws.Range("J6").Value = 30000
$ws.Range("K6").Value = 531
$ws.Range("L6").Value = 90000
$ws.Range("M6").Value = -419
$ws.Range("N6").Value = -90224
$ws.Range("H8").Value = 224.28572
$ws.Range("I8").Value = 224.28572
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 672.85716
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -533.85716
$ws.Range("H31").Value = 4331.8335
$ws.Range("I31").Value = 2998.4
$ws.Range("J31").Value = 10999
$ws.Range("K31").Value = 8995.200000000001
$ws.Range("L31").Value = 32997
$ws.Range("M31").Value = -8765.200000000001
$ws.Range("N31").Value = -33457
$ws.Range("H132").Value = 2541.027
$ws.Range("I132").Value = 2607.7188
$ws.Range("J132").Value = 2114.2
$ws.Range("K132").Value = 7823.1564
$ws.Range("L132").Value = 6342.599999999999
$ws.Range("M132").Value = -5293.1564
$ws.Range("N132").Value = -11402.6
$ws.Range("H137").Value = 2635.3845
$ws.Range("I137").Value = 1465.4
$ws.Range("J137").Value = 3366.625
$ws.Range("K137").Value = 4396.200000000001
$ws.Range("L137").Value = 10099.875
$ws.Range("M137").Value = -1846.200000000001
$ws.Range("N137").Value = -15199.875

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7698.1763
$ws.Range("I61").Value = 6765.5835
$ws.Range("J61").Value = 9936.4
$ws.Range("K61").Value = 6765.5835
$ws.Range("L61").Value = 9936.4
$ws.Range("M61").Value = -6553.5835
$ws.Range("N61").Value = -10360.4
$ws.Range("H74").Value = 1545.3871
$ws.Range("I74").Value = 1477.4166
$ws.Range("J74").Value = 1778.4286
$ws.Range("K74").Value = 1477.4166
$ws.Range("L74").Value = 1778.4286
$ws.Range("M74").Value = -603.4166
$ws.Range("N74").Value = -3526.4286
$ws.Range("H77").Value = 1545.3871
$ws.Range("I77").Value = 1477.4166
$ws.Range("J77").Value = 1778.4286
$ws.Range("K77").Value = 7387.083000000001
$ws.Range("L77").Value = 8892.143
$ws.Range("M77").Value = -3019.083000000001
$ws.Range("N77").Value = -17628.143
$ws.Range("H97").Value = 4999.5
$ws.Range("I97").Value = 4999.5
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 4999.5
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -4503.5
$ws.Range("N97").ClearContents()
$ws.Range("H132").Value = 2602.5483
$ws.Range("I132").Value = 2272.4075
$ws.Range("J132").Value = 4831
$ws.Range("K132").Value = 6817.2225
$ws.Range("L132").Value = 14493
$ws.Range("M132").Value = -4287.2225
$ws.Range("N132").Value = -19553
$ws.Range("H136").Value = 7698.1763
$ws.Range("I136").Value = 6765.5835
$ws.Range("J136").Value = 9936.4
$ws.Range("K136").Value = 20296.7505
$ws.Range("L136").Value = 29809.2
$ws.Range("M136").Value = -17746.7505
$ws.Range("N136").Value = -34909.2

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 445.27274
$ws.Range("I22").Value = 428.57144
$ws.Range("J22").Value = 474.5
$ws.Range("K22").Value = 428.57144
$ws.Range("L22").Value = 474.5
$ws.Range("M22").Value = -255.57144
$ws.Range("N22").Value = -820.5
$ws.Range("H105").Value = 1487.5
$ws.Range("I105").Value = 1487.5
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1487.5
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 259.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 66519.086
$ws.Range("I16").Value = 15823
$ws.Range("J16").Value = 319999.5
$ws.Range("K16").Value = 15823
$ws.Range("L16").Value = 319999.5
$ws.Range("M16").Value = -15536
$ws.Range("N16").Value = -320573.5
$ws.Range("H31").Value = 3392.1428
$ws.Range("I31").Value = 3324.0881
$ws.Range("J31").Value = 3546.4
$ws.Range("K31").Value = 3324.0881
$ws.Range("L31").Value = 3546.4
$ws.Range("M31").Value = -3029.0881
$ws.Range("N31").Value = -4136.4
$ws.Range("H34").Value = 3392.1428
$ws.Range("I34").Value = 3324.0881
$ws.Range("J34").Value = 3546.4
$ws.Range("K34").Value = 3324.0881
$ws.Range("L34").Value = 3546.4
$ws.Range("M34").Value = -3122.0881
$ws.Range("N34").Value = -3950.4
$ws.Range("H113").Value = 66519.086
$ws.Range("I113").Value = 15823
$ws.Range("J113").Value = 319999.5
$ws.Range("K113").Value = 15823
$ws.Range("L113").Value = 319999.5
$ws.Range("M113").Value = -13653
$ws.Range("N113").Value = -324339.5
$ws.Range("H132").Value = 3218.4583
$ws.Range("I132").Value = 2642.476
$ws.Range("J132").Value = 7250.3335
$ws.Range("K132").Value = 7927.428
$ws.Range("L132").Value = 21751.0005
$ws.Range("M132").Value = -5397.428
$ws.Range("N132").Value = -26811.0005
$ws.Range("H134").Value = 5566.59
$ws.Range("I134").Value = 4795.9414
$ws.Range("J134").Value = 10807
$ws.Range("K134").Value = 14387.8242
$ws.Range("L134").Value = 32421
$ws.Range("M134").Value = -11852.8242
$ws.Range("N134").Value = -37491

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 77287860
$ws.Range("I4").Value = 100069210
$ws.Range("J4").Value = 1350000
$ws.Range("K4").Value = 300207630
$ws.Range("L4").Value = 4050000
$ws.Range("M4").Value = -300207518
$ws.Range("N4").Value = -4050224
$ws.Range("H6").Value = 97.55556
$ws.Range("I6").Value = 104.75
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 314.25
$ws.Range("L6").Value = 120
$ws.Range("M6").Value = -201.25
$ws.Range("N6").Value = -346
$ws.Range("H12").Value = 298.65714
$ws.Range("I12").Value = 285.58334
$ws.Range("J12").Value = 305.47827
$ws.Range("K12").Value = 856.7500200000001
$ws.Range("L12").Value = 916.43481
$ws.Range("M12").Value = -683.7500200000001
$ws.Range("N12").Value = -1262.43481
$ws.Range("H81").Value = 3249.25
$ws.Range("I81").Value = 3000
$ws.Range("J81").Value = 3498.5
$ws.Range("K81").Value = 9000
$ws.Range("L81").Value = 10495.5
$ws.Range("M81").Value = -7877
$ws.Range("N81").Value = -12741.5
$ws.Range("H84").Value = 3249.25
$ws.Range("I84").Value = 3000
$ws.Range("J84").Value = 3498.5
$ws.Range("K84").Value = 27000
$ws.Range("L84").Value = 31486.5
$ws.Range("M84").Value = -21384
$ws.Range("N84").Value = -42718.5
$ws.Range("H98").Value = 1869.25
$ws.Range("I98").Value = 1192.6
$ws.Range("J98").Value = 2997
$ws.Range("K98").Value = 3577.8
$ws.Range("L98").Value = 8991
$ws.Range("M98").Value = -2079.8
$ws.Range("N98").Value = -11987

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 59999
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 59999
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 59999
$ws.Range("N68").Value = -61621
$ws.Range("H70").Value = 6652.5293
$ws.Range("I70").Value = 6672.8667
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 6672.8667
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -6402.8667
$ws.Range("N70").Value = -7040
$ws.Range("H71").Value = 59999
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 59999
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 179997
$ws.Range("N71").Value = -188109
$ws.Range("H73").Value = 6652.5293
$ws.Range("I73").Value = 6672.8667
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 6672.8667
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = -5736.8667
$ws.Range("N73").Value = -8372
$ws.Range("H80").Value = 3457.375
$ws.Range("I80").Value = 1234.6666
$ws.Range("J80").Value = 4791
$ws.Range("K80").Value = 1234.6666
$ws.Range("L80").Value = 4791
$ws.Range("M80").Value = -236.6666
$ws.Range("N80").Value = -6787
$ws.Range("H83").Value = 3457.375
$ws.Range("I83").Value = 1234.6666
$ws.Range("J83").Value = 4791
$ws.Range("K83").Value = 6173.333000000001
$ws.Range("L83").Value = 23955
$ws.Range("M83").Value = -1181.333000000001
$ws.Range("N83").Value = -33939

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2775.5454
$ws.Range("I22").Value = 2600.1428
$ws.Range("J22").Value = 3082.5
$ws.Range("K22").Value = 2600.1428
$ws.Range("L22").Value = 3082.5
$ws.Range("M22").Value = -2305.1428
$ws.Range("N22").Value = -3672.5
$ws.Range("H27").Value = 2775.5454
$ws.Range("I27").Value = 2600.1428
$ws.Range("J27").Value = 3082.5
$ws.Range("K27").Value = 2600.1428
$ws.Range("L27").Value = 3082.5
$ws.Range("M27").Value = -2493.1428
$ws.Range("N27").Value = -3296.5
$ws.Range("H55").Value = 219.53572
$ws.Range("I55").Value = 258.9375
$ws.Range("J55").Value = 167
$ws.Range("K55").Value = 258.9375
$ws.Range("L55").Value = 167
$ws.Range("M55").Value = -85.9375
$ws.Range("N55").Value = -513
$ws.Range("H82").Value = 1614.4
$ws.Range("I82").Value = 1083.3334
$ws.Range("J82").Value = 1842
$ws.Range("K82").Value = 1083.3334
$ws.Range("L82").Value = 1842
$ws.Range("M82").Value = -722.3334
$ws.Range("N82").Value = -2564
$ws.Range("H85").Value = 1614.4
$ws.Range("I85").Value = 1083.3334
$ws.Range("J85").Value = 1842
$ws.Range("K85").Value = 1083.3334
$ws.Range("L85").Value = 1842
$ws.Range("M85").Value = 164.6666
$ws.Range("N85").Value = -4338
$ws.Range("H93").Value = 2291.1667
$ws.Range("I93").Value = 1349.4
$ws.Range("J93").Value = 7000
$ws.Range("K93").Value = 1349.4
$ws.Range("L93").Value = 7000
$ws.Range("M93").Value = -101.4000000000001
$ws.Range("N93").Value = -9496
$ws.Range("H100").Value = 2000.75
$ws.Range("I100").Value = 2000.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2000.75
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1459.75
$ws.Range("H133").Value = 84332.336
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 84332.336
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 84332.336
$ws.Range("N133").Value = -89392.336

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H132").Value = 2141.422
$ws.Range("I132").Value = 2312.0667
$ws.Range("J132").Value = 1800.1333
$ws.Range("K132").Value = 6936.2001
$ws.Range("L132").Value = 5400.3999
$ws.Range("M132").Value = -4406.2001
$ws.Range("N132").Value = -10460.3999

Write-Host "Applied all cell updates"